# feat: Add direction handling for divergences in cruzar function
#
# Adds a new "dir" column (E) to the "divergencias" sheet: a header cell
# styled like the existing header row, a sample data value ("MENOR") for
# the first data row, and a matching column width.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("divergencias")

# New header cell E1: copy formatting from the neighboring header (D1,
# "dif_rel") so it matches the bold/centered/bordered header style, then
# set its own text.
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value = "dir"

# New data cell E2 (plain/default style, like the other text cells).
$ws.Range("E2").Value = "MENOR"

# New column width (stored width 7 == ColumnWidth 6.17 once Excel's
# character-padding offset is accounted for).
$ws.Columns.Item(5).ColumnWidth = 6.17
